$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the product name text (shared string used by A2)
$ws.Range("A2").Value = "Оффлайн ТВ 2 плитки"

# Update price and quantity values
$ws.Range("B2").Value = 2300
$ws.Range("C2").Value = 160

# Update selection to match the saved view state (active cell C2, single cell selected)
$ws.Range("C2").Select()
